# The "10171" sheet (Codelists - Import / currency codelist) was missing
# the ISO 4217 "KMF" (Comorian Franc) entry that all the other rows use
# as their code key; only a malformed " KMF" (leading-space) row existed
# at row 78. Re-insert a proper "KMF" row right after it, at row 79,
# pushing every following currency row (KPW..ZWL) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10171")

$firstShiftRow = 79
$lastOriginalRow = 181

# Walk bottom-up so each row is copied to its new home before it gets
# overwritten by the row above sliding down into it.
for ($r = $lastOriginalRow; $r -ge $firstShiftRow; $r--) {
    $src = $ws.Range("A" + $r + ":B" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":B" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# Row 79 is now a duplicate of row 80 (KPW) — overwrite it with the
# new KMF entry the edit introduces. Code and description both read
# "KMF" (matching the source data exactly).
$ws.Range("A79").Value2 = "KMF"
$ws.Range("B79").Value2 = "KMF"

# Row 182 is a brand-new row with no banding style of its own yet
# (rows 79-181 inherited correct banding automatically because we only
# overwrote Value2, never their formatting). Clone the zebra-stripe
# formatting from row 180 — same (even) parity — onto row 182 so
# columns A:D keep the alternating band look the rest of the table has.
$ws.Range("A180:D180").Copy()
$ws.Range("A182:D182").PasteSpecial(-4122)
